# Add a new "everyone wins" tournament machine worksheet, modeled after the
# existing "tournament3" sheet (same layout/columns/number-format styling),
# but with Cherry/Heart/Orange/Gold Bar/Seven symbols and a payout table
# where every combination pays out (no "loses" distinction).

$wb = $excel.ActiveWorkbook

# Duplicate "tournament3" so the new sheet inherits identical column widths,
# number formats, and page setup, then place it right after tournament3.
$src = $wb.Sheets("tournament3")
$src.Copy([Type]::Missing, $src) | Out-Null
$ws = $wb.Sheets($src.Index + 1)
$ws.Name = "tournament4"

# Wipe the copied content -- we'll rebuild it from scratch for the new machine.
$ws.Cells.Clear() | Out-Null

# --- Header row for the symbol table ---
$ws.Range("B1").Value = "Slot 1"
$ws.Range("C1").Value = "Slot 2"
$ws.Range("D1").Value = "Slot 3"

# --- Symbols + per-reel counts ---
$ws.Range("A2").Value = "Cherry"
$ws.Range("B2").Value = 6
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = 8

$ws.Range("A3").Value = "Heart"
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 6
$ws.Range("D3").Value = 4

$ws.Range("A4").Value = "Orange"
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 5

$ws.Range("A5").Value = "Gold Bar"
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 3

$ws.Range("A6").Value = "Seven"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 4

# --- Total combinations ---
$ws.Range("A8").Value = "Combinations"
$ws.Range("B8").Formula = "=SUM(B2:B7)*SUM(C2:C7)*SUM(D2:D7)"

# --- Payout table header ---
$ws.Range("B10").Value = "Number"
$ws.Range("C10").Value = "Frequency"
$ws.Range("D10").Value = "Payout"
$ws.Range("E10").Value = "Payout Ratio"

# --- Cherry ---
$ws.Range("A11").Value = "Cherry+Any+Any"
$ws.Range("B11").Formula = "=B2*SUM(C2:C6)*SUM(D2:D6)-B12"
$ws.Range("C11").Formula = '=B11/$B$8'
$ws.Range("D11").Value = 0.2
$ws.Range("E11").Formula = "=C11*D11"

$ws.Range("A12").Value = "Cherry+Cherry+Cherry"
$ws.Range("B12").Formula = "=B2*C2*D2"
$ws.Range("C12").Formula = '=B12/$B$8'
$ws.Range("D12").Value = 4
$ws.Range("E12").Formula = "=C12*D12"

# --- Heart ---
$ws.Range("A13").Value = "Heart+Any+Any"
$ws.Range("B13").Formula = "=B3*SUM(C2:C6)*SUM(D2:D6)-B14"
$ws.Range("C13").Formula = '=B13/$B$8'
$ws.Range("D13").Value = 0.2
$ws.Range("E13").Formula = "=C13*D13"

$ws.Range("A14").Value = "Heart+Heart+Heart"
$ws.Range("B14").Formula = "=B3*C3*D3"
$ws.Range("C14").Formula = '=B14/$B$8'
$ws.Range("D14").Value = 6
$ws.Range("E14").Formula = "=C14*D14"

# --- Orange ---
$ws.Range("A15").Value = "Orange+Any+Any"
$ws.Range("B15").Formula = "=B4*SUM(C2:C6)*SUM(D2:D6)-B16"
$ws.Range("C15").Formula = '=B15/$B$8'
$ws.Range("D15").Value = 0.2
$ws.Range("E15").Formula = "=C15*D15"

$ws.Range("A16").Value = "Orange+Orange+Orange"
$ws.Range("B16").Formula = "=B4*C4*D4"
$ws.Range("C16").Formula = '=B16/$B$8'
$ws.Range("D16").Value = 8
$ws.Range("E16").Formula = "=C16*D16"

# --- Gold Bar ---
$ws.Range("A17").Value = "Gold Bar+Any+Any"
$ws.Range("B17").Formula = "=B5*SUM(C2:C6)*SUM(D2:D6)-B18"
$ws.Range("C17").Formula = '=B17/$B$8'
$ws.Range("D17").Value = 0.2
$ws.Range("E17").Formula = "=C17*D17"

$ws.Range("A18").Value = "Gold Bar+Gold Bar+Gold Bar"
$ws.Range("B18").Formula = "=B5*C5*D5"
$ws.Range("C18").Formula = '=B18/$B$8'
$ws.Range("D18").Value = 10
$ws.Range("E18").Formula = "=C18*D18"

# --- Seven ---
$ws.Range("A19").Value = "Seven+Any+Any"
$ws.Range("B19").Formula = "=B6*SUM(C2:C5)*SUM(D2:D6)"
$ws.Range("C19").Formula = '=B19/$B$8'
$ws.Range("D19").Value = 10
$ws.Range("E19").Formula = "=C19*D19"

$ws.Range("A20").Value = "Seven+Seven+Any"
$ws.Range("B20").Formula = "=B6*C6*SUM(D2:D5)"
$ws.Range("C20").Formula = '=B20/$B$8'
$ws.Range("D20").Value = 20
$ws.Range("E20").Formula = "=C20*D20"

$ws.Range("A21").Value = "Seven+Seven+Seven"
$ws.Range("B21").Formula = "=B6*C6*D6"
$ws.Range("C21").Formula = '=B21/$B$8'
$ws.Range("D21").Value = 50
$ws.Range("E21").Formula = "=C21*D21"

# Blank spacer row, keeps the ratio-column formatting like the source sheet.
$ws.Range("E22").NumberFormat = "0.000"

# --- Totals ---
$ws.Range("A23").Value = "Total"
$ws.Range("C23").Formula = "=SUM(C11:C22)"
$ws.Range("E23").Formula = "=SUM(E11:E22)"

$ws.Range("C24").Formula = '=SUMIF(D11:D21,">1",C11:C21)'

# Match the source sheet's number formats for the ratio/payout columns.
$ws.Range("C11:C21").NumberFormat = "0.0000"
$ws.Range("E11:E21").NumberFormat = "0.000"
$ws.Range("C23:C24").NumberFormat = "0.000"
$ws.Range("E23").NumberFormat = "0.000"

# Make the new sheet the active tab with A7 selected, matching the edit.
$ws.Activate()
$ws.Range("A7").Select() | Out-Null
